$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.062.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "'3.773.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'629.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'165.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'3.771.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'6.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("D14").Value = "'34.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").Value = "'4.409.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "'3.783.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "'69.077.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'17.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.54%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'7.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "'468.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'9.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").Value = "'0.703"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'82.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("E25").Value = "  -7.43%  "
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "'10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D30").Value = "'3.922.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("D32").Value = "'2.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "'7.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").Value = "'0.179"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +20.25%  "
$ws.Range("D35").Value = "'28.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'3.725.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'8.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").Value = "'3.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.69%  "
$ws.Range("D41").Value = "'5.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D43").Value = "'0.962"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'2.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.83%  "
$ws.Range("D46").Value = "'156.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'43.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.93%  "
$ws.Range("D48").Value = "'47.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("D50").Value = "'0.294"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("D51").Value = "'8.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.98%  "

Write-Output "Updated cryptos list"
